$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tgfb2"
$ws.Range("C2").Value = "Tgfbr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.787349
$ws.Range("H2").Value = 5.362047
$ws.Range("I2").Value = 0.04925512201701282
$ws.Range("J2").Value = 0.04925512201701282
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 19.86261233333333
$ws.Range("N2").Value = 59.587837
$ws.Range("O2").Value = 0.1710751304955294
$ws.Range("P2").Value = 0.1710751304955294
$ws.Range("Q2").Value = 35.50142029137101
$ws.Range("R2").Value = 319.512782622339
$ws.Range("S2").Value = 0.008426326426633693
$ws.Range("T2").Value = 0.008426326426633693

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tgfb2"
$ws.Range("C3").Value = "Tgfbr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.787349
$ws.Range("H3").Value = 5.362047
$ws.Range("I3").Value = 0.04925512201701282
$ws.Range("J3").Value = 0.04925512201701282
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 69.67747766666666
$ws.Range("N3").Value = 209.032433
$ws.Range("O3").Value = 0.6001266794307873
$ws.Range("P3").Value = 0.6001266794307873
$ws.Range("Q3").Value = 124.537970030039
$ws.Range("R3").Value = 1120.841730270351
$ws.Range("S3").Value = 0.02955931282102817
$ws.Range("T3").Value = 0.02955931282102817

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tgfb2"
$ws.Range("C4").Value = "Tgfbr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.787349
$ws.Range("H4").Value = 5.362047
$ws.Range("I4").Value = 0.04925512201701282
$ws.Range("J4").Value = 0.04925512201701282
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.564526
$ws.Range("N4").Value = 79.693578
$ws.Range("O4").Value = 0.2287981900736832
$ws.Range("P4").Value = 0.2287981900736832
$ws.Range("Q4").Value = 47.480078981574
$ws.Range("R4").Value = 427.3207108341661
$ws.Range("S4").Value = 0.01126948276935096
$ws.Range("T4").Value = 0.01126948276935096

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tgfb2"
$ws.Range("C5").Value = "Tgfbr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 19.46983
$ws.Range("H5").Value = 58.40949000000001
$ws.Range("I5").Value = 0.5365425847444997
$ws.Range("J5").Value = 0.5365425847444997
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 19.86261233333333
$ws.Range("N5").Value = 59.587837
$ws.Range("O5").Value = 0.1710751304955294
$ws.Range("P5").Value = 0.1710751304955294
$ws.Range("Q5").Value = 386.7216854859034
$ws.Range("R5").Value = 3480.49516937313
$ws.Range("S5").Value = 0.09178909270157394
$ws.Range("T5").Value = 0.09178909270157394

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tgfb2"
$ws.Range("C6").Value = "Tgfbr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 19.46983
$ws.Range("H6").Value = 58.40949000000001
$ws.Range("I6").Value = 0.5365425847444997
$ws.Range("J6").Value = 0.5365425847444997
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 69.67747766666666
$ws.Range("N6").Value = 209.032433
$ws.Range("O6").Value = 0.6001266794307873
$ws.Range("P6").Value = 0.6001266794307873
$ws.Range("Q6").Value = 1356.608644998797
$ws.Range("R6").Value = 12209.47780498917
$ws.Range("S6").Value = 0.3219935197559284
$ws.Range("T6").Value = 0.3219935197559284

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tgfb2"
$ws.Range("C7").Value = "Tgfbr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 19.46983
$ws.Range("H7").Value = 58.40949000000001
$ws.Range("I7").Value = 0.5365425847444997
$ws.Range("J7").Value = 0.5365425847444997
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 26.564526
$ws.Range("N7").Value = 79.693578
$ws.Range("O7").Value = 0.2287981900736832
$ws.Range("P7").Value = 0.2287981900736832
$ws.Range("Q7").Value = 517.20680525058
$ws.Range("R7").Value = 4654.86124725522
$ws.Range("S7").Value = 0.1227599722869973
$ws.Range("T7").Value = 0.1227599722869973

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Tgfb2"
$ws.Range("C8").Value = "Tgfbr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 15.03039733333333
$ws.Range("H8").Value = 45.091192
$ws.Range("I8").Value = 0.4142022932384875
$ws.Range("J8").Value = 0.4142022932384875
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 19.86261233333333
$ws.Range("N8").Value = 59.587837
$ws.Range("O8").Value = 0.1710751304955294
$ws.Range("P8").Value = 0.1710751304955294
$ws.Range("Q8").Value = 298.5429554479671
$ws.Range("R8").Value = 2686.886599031704
$ws.Range("S8").Value = 0.0708597113673218
$ws.Range("T8").Value = 0.0708597113673218

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Tgfb2"
$ws.Range("C9").Value = "Tgfbr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 15.03039733333333
$ws.Range("H9").Value = 45.091192
$ws.Range("I9").Value = 0.4142022932384875
$ws.Range("J9").Value = 0.4142022932384875
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 69.67747766666666
$ws.Range("N9").Value = 209.032433
$ws.Range("O9").Value = 0.6001266794307873
$ws.Range("P9").Value = 0.6001266794307873
$ws.Range("Q9").Value = 1047.28017451446
$ws.Range("R9").Value = 9425.521570630135
$ws.Range("S9").Value = 0.2485738468538308
$ws.Range("T9").Value = 0.2485738468538308

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tgfb2"
$ws.Range("C10").Value = "Tgfbr3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 15.03039733333333
$ws.Range("H10").Value = 45.091192
$ws.Range("I10").Value = 0.4142022932384875
$ws.Range("J10").Value = 0.4142022932384875
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 26.564526
$ws.Range("N10").Value = 79.693578
$ws.Range("O10").Value = 0.2287981900736832
$ws.Range("P10").Value = 0.2287981900736832
$ws.Range("Q10").Value = 399.275380751664
$ws.Range("R10").Value = 3593.478426764976
$ws.Range("S10").Value = 0.09476873501733493
$ws.Range("T10").Value = 0.09476873501733495
